$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts columns B:F left to A:E,
# and removes the bold-bordered style that was on A2:A3.
$ws.Range("A:A").Delete()
